$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D sometimes holds plain numeric-looking text (e.g. "0.699",
# "103.68", "0.0000104") that must stay TEXT, matching the source
# workbook (all data cells are stored as inline/shared strings, never
# numbers). Assigning such strings directly would let Excel silently
# coerce them into real numbers (and even scientific notation for very
# small values), so we temporarily force Text format on the whole
# price column, make the assignments, then clear the format override
# again so cell styling matches the original (unstyled) cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Formula = '44.229.52'
$ws.Range('E2').Formula = '  +2.21%  '
$ws.Range('D3').Formula = '2.386.29'
$ws.Range('E3').Formula = '  +1.32%  '
$ws.Range('E4').Formula = '  -0.05%  '
$ws.Range('D5').Formula = '0.699'
$ws.Range('E5').Formula = '  +8.03%  '
$ws.Range('D6').Formula = '244.04'
$ws.Range('E6').Formula = '  +4.84%  '
$ws.Range('D7').Formula = '77.31'
$ws.Range('E7').Formula = '  +7.40%  '
$ws.Range('E8').Formula = '  -0.10%  '
$ws.Range('D9').Formula = '0.631'
$ws.Range('E9').Formula = '  +32.24%  '
$ws.Range('D10').Formula = '0.105'
$ws.Range('E10').Formula = '  +6.96%  '
$ws.Range('D11').Formula = '58.12'
$ws.Range('E11').Formula = '  +2.37%  '
$ws.Range('D12').Formula = '33.63'
$ws.Range('E12').Formula = '  +24.83%  '
$ws.Range('D13').Formula = '7.63'
$ws.Range('E13').Formula = '  +21.87%  '
$ws.Range('E14').Formula = '  +2.34%  '
$ws.Range('B15').Formula = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Formula = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Formula = '2.741.18'
$ws.Range('E15').Formula = '  +1.27%  '
$ws.Range('B16').Formula = 'Chainlink'
$ws.Range('C16').Formula = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Formula = '17.28'
$ws.Range('E16').Formula = '  +7.90%  '
$ws.Range('D17').Formula = '0.935'
$ws.Range('E17').Formula = '  +8.45%  '
$ws.Range('D18').Formula = '2.385.65'
$ws.Range('E18').Formula = '  +1.32%  '
$ws.Range('D19').Formula = '44.236.21'
$ws.Range('E19').Formula = '  +2.30%  '
$ws.Range('D20').Formula = '0.0000104'
$ws.Range('E20').Formula = '  +2.73%  '
$ws.Range('E21').Formula = '  +7.08%  '
$ws.Range('D22').Formula = '79.19'
$ws.Range('E22').Formula = '  +6.73%  '
$ws.Range('D23').Formula = '259.67'
$ws.Range('E23').Formula = '  +4.09%  '
$ws.Range('E24').Formula = '  +0.12%  '
$ws.Range('E25').Formula = '  +5.03%  '
$ws.Range('D26').Formula = '3.72'
$ws.Range('E26').Formula = '  +1.08%  '
$ws.Range('D27').Formula = '11.05'
$ws.Range('E27').Formula = '  +10.73%  '
$ws.Range('E28').Formula = '  +19.46%  '
$ws.Range('E29').Formula = '  +1.79%  '
$ws.Range('D30').Formula = '23.26'
$ws.Range('E30').Formula = '  +4.30%  '
$ws.Range('D31').Formula = '175.65'
$ws.Range('E31').Formula = '  +1.08%  '
$ws.Range('E32').Formula = '  +2.48%  '
$ws.Range('D33').Formula = '0.137'
$ws.Range('E33').Formula = '  +7.89%  '
$ws.Range('E34').Formula = '  +9.67%  '
$ws.Range('D35').Formula = '0.0767'
$ws.Range('E35').Formula = '  +11.15%  '
$ws.Range('D36').Formula = '5.41'
$ws.Range('E36').Formula = '  +7.89%  '
$ws.Range('E37').Formula = '  +6.50%  '
$ws.Range('E38').Formula = '  +3.36%  '
$ws.Range('D39').Formula = '6.64'
$ws.Range('E39').Formula = '  +1.18%  '
$ws.Range('E40').Formula = '  +9.92%  '
$ws.Range('D41').Formula = '9.20'
$ws.Range('E41').Formula = '  +3.61%  '
$ws.Range('D42').Formula = '19.23'
$ws.Range('E42').Formula = '  +4.64%  '
$ws.Range('B43').Formula = 'Algorand'
$ws.Range('C43').Formula = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Formula = '0.201'
$ws.Range('E43').Formula = '  +20.51%  '
$ws.Range('B44').Formula = 'BinanceUSD'
$ws.Range('C44').Formula = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D44').Formula = '1.00'
$ws.Range('E44').Formula = '  -0.06%  '
$ws.Range('E45').Formula = '  +5.92%  '
$ws.Range('D46').Formula = '2.57'
$ws.Range('E46').Formula = '  +16.28%  '
$ws.Range('D47').Formula = '1.29'
$ws.Range('E47').Formula = '  +6.47%  '
$ws.Range('B48').Formula = 'Aave'
$ws.Range('C48').Formula = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Formula = '103.68'
$ws.Range('E48').Formula = '  +3.71%  '
$ws.Range('B49').Formula = 'Cronos'
$ws.Range('C49').Formula = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Formula = '0.101'
$ws.Range('E49').Formula = '  +6.30%  '
$ws.Range('D50').Formula = '4.64'
$ws.Range('E50').Formula = '  +3.10%  '
$ws.Range('D51').Formula = '55.39'
$ws.Range('E51').Formula = '  +10.43%  '

$ws.Range("D2:D51").ClearFormats()
